$d = $word.ActiveDocument

# Locate the paragraph ending with "Text parsed from text boxes and tables."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Text parsed from text boxes and tables.*") {
        $target = $p
        break
    }
}

# Insert a new paragraph right after it, inheriting formatting, then
# set its text and promote it one list level up (ilvl 2 -> 1) so it
# becomes a sibling of "Slide Name" / "Slide Number" / "Slide Text"
# rather than a child of "Slide Text".
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.Text = "A count of the number of shapes on the slide. "
$newPara.Range.ListFormat.ListLevelNumber = 2
